$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 37: section header for the new sensor-location experiment ---
$ws.Range("F37").Value = "new sensor location experiment"

# --- Row 38: existing "leave 1 out waist" label (A38) stays; add the new table header row ---
$ws.Range("F38").Value = "400, 50, maxmag"
$ws.Range("G38").Value = "waist"
$ws.Range("H38").Value = "left"
$ws.Range("I38").Value = "right"
$ws.Range("J38").Value = "W L"
$ws.Range("K38").Value = "W R"
$ws.Range("L38").Value = "R L"
$ws.Range("M38").Value = "all"

# --- Row 39: existing "full f1 avg"/0.82 (A39/B39) stay; add the "full acc" data row ---
$ws.Range("F39").Value = "full acc"
$ws.Range("G39").Value = 0.81916067969252504
$ws.Range("H39").Value = 0.85101199171269004
$ws.Range("I39").Value = 0.81825354459291699
$ws.Range("J39").Value = 0.88158911832264197
$ws.Range("K39").Value = 0.87687349263494896
$ws.Range("L39").Value = 0.91159276762340302
$ws.Range("M39").Value = 0.916500755539568

# --- New row 40: "jump f1" data row ---
$ws.Range("F40").Value = "jump f1"
$ws.Range("G40").Value = 0.85775351128751598
$ws.Range("H40").Value = 0.90739485123363395
$ws.Range("I40").Value = 0.88217792117237903
$ws.Range("J40").Value = 0.91170816872625104
$ws.Range("K40").Value = 0.92854280278030998
$ws.Range("L40").Value = 0.95385457966846099
$ws.Range("M40").Value = 0.95285580712696105

# --- Move/resize the "bin size experiment" chart (Chart 5) down to make room ---
$co = $ws.ChartObjects(4)
$co.Top = 87.6
$co.Left = 294.298828125

# --- Update the view: scroll position + selection ---
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("M40").Select()
